$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the original column A (shifts B:F left to A:E)
$ws.Range("A1:A3").Delete([Microsoft.Office.Interop.Excel.XlDeleteShiftDirection]::xlShiftToLeft)
